$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the cell's original
# (default/no-op) style - Excel's Value-assignment auto-detects numeric-
# looking strings (e.g. "1.00", "0.999") and would otherwise silently
# convert them to numbers, which is not what the source data represents.
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "51.543.99"
Set-TextValue $ws.Range("E2") "  +0.65%  "
Set-TextValue $ws.Range("D3") "2.980.47"
Set-TextValue $ws.Range("E3") "  +2.23%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.02%  "
Set-TextValue $ws.Range("D5") "379.30"
Set-TextValue $ws.Range("E5") "  +2.58%  "
Set-TextValue $ws.Range("D6") "104.82"
Set-TextValue $ws.Range("E6") "  +0.86%  "
Set-TextValue $ws.Range("E7") "  +0.16%  "
Set-TextValue $ws.Range("D8") "1.00"
Set-TextValue $ws.Range("E8") "  +0.02%  "
Set-TextValue $ws.Range("D9") "0.594"
Set-TextValue $ws.Range("E9") "  +1.38%  "
Set-TextValue $ws.Range("D10") "37.23"
Set-TextValue $ws.Range("E10") "  +1.57%  "
Set-TextValue $ws.Range("E11") "  -0.10%  "
Set-TextValue $ws.Range("D12") "0.0845"
Set-TextValue $ws.Range("E12") "  +1.27%  "
Set-TextValue $ws.Range("D13") "3.449.83"
Set-TextValue $ws.Range("E13") "  +2.16%  "
Set-TextValue $ws.Range("D14") "18.43"
Set-TextValue $ws.Range("E14") "  +0.43%  "
Set-TextValue $ws.Range("D15") "7.59"
Set-TextValue $ws.Range("E15") "  +2.65%  "
Set-TextValue $ws.Range("D16") "2.973.94"
Set-TextValue $ws.Range("E16") "  +1.79%  "
Set-TextValue $ws.Range("D17") "0.970"
Set-TextValue $ws.Range("E17") "  +3.91%  "
Set-TextValue $ws.Range("D18") "51.507.61"
Set-TextValue $ws.Range("E18") "  +0.66%  "
Set-TextValue $ws.Range("E19") "  +2.66%  "
Set-TextValue $ws.Range("E20") "  +3.20%  "
Set-TextValue $ws.Range("D21") "12.96"
Set-TextValue $ws.Range("E21") "  +0.08%  "
Set-TextValue $ws.Range("E22") "  +2.01%  "
Set-TextValue $ws.Range("D23") "69.38"
Set-TextValue $ws.Range("E23") "  +1.39%  "
Set-TextValue $ws.Range("D24") "262.10"
Set-TextValue $ws.Range("E24") "  +0.97%  "
Set-TextValue $ws.Range("D25") "2.83"
Set-TextValue $ws.Range("E25") "  +5.63%  "
Set-TextValue $ws.Range("D26") "8.20"
Set-TextValue $ws.Range("E26") "  +16.11%  "
Set-TextValue $ws.Range("D27") "7.69"
Set-TextValue $ws.Range("E27") "  +23.82%  "
Set-TextValue $ws.Range("D28") "0.116"
Set-TextValue $ws.Range("E28") "  +12.27%  "
Set-TextValue $ws.Range("D29") "0.171"
Set-TextValue $ws.Range("E29") "  -2.83%  "
Set-TextValue $ws.Range("E30") "  +0.11%  "
Set-TextValue $ws.Range("D31") "25.90"
Set-TextValue $ws.Range("E31") "  +0.65%  "
Set-TextValue $ws.Range("D32") "9.85"
Set-TextValue $ws.Range("E32") "  -0.37%  "
Set-TextValue $ws.Range("D33") "35.06"
Set-TextValue $ws.Range("E33") "  +1.27%  "
Set-TextValue $ws.Range("E34") "  -2.08%  "
Set-TextValue $ws.Range("D35") "51.10"
Set-TextValue $ws.Range("E35") "  +0.52%  "
Set-TextValue $ws.Range("E36") "  +5.09%  "
Set-TextValue $ws.Range("E37") "  -0.02%  "
Set-TextValue $ws.Range("E38") "  +0.43%  "
Set-TextValue $ws.Range("D39") "17.20"
Set-TextValue $ws.Range("E39") "  +0.58%  "
Set-TextValue $ws.Range("E40") "  -2.54%  "
Set-TextValue $ws.Range("E41") "  +0.51%  "
Set-TextValue $ws.Range("D42") "0.115"
Set-TextValue $ws.Range("E42") "  +2.32%  "
Set-TextValue $ws.Range("D43") "125.55"
Set-TextValue $ws.Range("E43") "  +5.96%  "
Set-TextValue $ws.Range("D44") "21.60"
Set-TextValue $ws.Range("E44") "  -2.64%  "
Set-TextValue $ws.Range("D45") "0.284"
Set-TextValue $ws.Range("E45") "  +17.88%  "
Set-TextValue $ws.Range("D46") "2.05"
Set-TextValue $ws.Range("E46") "  -1.27%  "
Set-TextValue $ws.Range("E47") "  +3.54%  "
Set-TextValue $ws.Range("D48") "2.035.66"
Set-TextValue $ws.Range("E48") "  +0.84%  "
Set-TextValue $ws.Range("E49") "  +2.21%  "
Set-TextValue $ws.Range("D50") "0.0335"
Set-TextValue $ws.Range("E50") "  +8.18%  "
Set-TextValue $ws.Range("D51") "58.45"
Set-TextValue $ws.Range("E51") "  +2.98%  "
